$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts old rows 13-24 down to 14-25, and heights realign automatically)
$ws.Rows.Item(13).Insert()

# The inserted row 13 picks up a stray styled-but-empty A13 cell; remove it entirely
$ws.Range("A13").Clear()

# Row 1
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

# Row 2
$ws.Range("B2").Value = "LOB1041"
$ws.Range("C2").Value = "LOB1041"

# Row 3
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Física Experimental II"
$ws.Range("C3").Value = " Física Experimental II"

# Row 4
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Experimental Physics II"
$ws.Range("C4").Value = "Experimental Physics II"

# Row 5
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"

# Row 6
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

# Row 7
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Row 8
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2018"
$ws.Range("C8").Value = "01/01/2018"

# Row 9
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EF-2,EM-2,EA-3,EP-3"
$ws.Range("C9").Value = "EF-2,EM-2,EA-3,EP-3"

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Verificação experimental dos conceitos básicos de hidrostática, hidrodinâmica, termodinâmica e ondas."
$ws.Range("C10").Value = "Verificação experimental dos conceitos básicos de hidrostática, hidrodinâmica, termodinâmica e ondas."

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Experimental verification of the basic concepts of hydrostatic, hydrodynamic, thermodynamic and waves."
$ws.Range("C11").Value = "Experimental verification of the basic concepts of hydrostatic, hydrodynamic, thermodynamic and waves."

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13
$ws.Range("B13").Value = "5817535 - Lucas Barboza Sarno da Silva"
$ws.Range("C13").Value = "5817535 - Lucas Barboza Sarno da Silva"

# Row 14
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Abordagem experimental de conceitos relacionados à mecânica dos fluidos, termodinâmica, oscilações e ondas."
$ws.Range("C14").Value = "Abordagem experimental de conceitos relacionados à mecânica dos fluidos, termodinâmica, oscilações e ondas."

# Row 15
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Experimental approach to concepts related to fluid mechanics, thermodynamics, oscillations and waves."
$ws.Range("C15").Value = "Experimental approach to concepts related to fluid mechanics, thermodynamics, oscillations and waves."

# Row 16
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1) Princípio de Stevin e Pascal2) Empuxo e Princípio de Arquimedes3) Tensão superficial4) Queda em um meio viscoso5) Sistema massa-mola6) Ondas mecânicas7) Calor, temperatura e capacidade do corpo de armazenar energia8) Dilatação linear9) Os meios de propagação de calor10) Calor específico e calor latente11) A lei de Boyle-Mariotte"
$ws.Range("C16").Value = "1) Princípio de Stevin e Pascal2) Empuxo e Princípio de Arquimedes3) Tensão superficial4) Queda em um meio viscoso5) Sistema massa-mola6) Ondas mecânicas7) Calor, temperatura e capacidade do corpo de armazenar energia8) Dilatação linear9) Os meios de propagação de calor10) Calor específico e calor latente11) A lei de Boyle-Mariotte"

# Row 17
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "1) stevin’s and Pascal’s Principle2) Thrust and Archimedes’ Principle3) Surface tension4) The fall in a viscous fluid5) Mass-spring system6) Mechanical waves7) Heat, temperature, and the body's capacity to store energy8) Linear thermal expansion9) The fundamental modes of heat transfer10) Specific and latent heat11) The Boyle-Mariotte’s Law"
$ws.Range("C17").Value = "1) stevin’s and Pascal’s Principle2) Thrust and Archimedes’ Principle3) Surface tension4) The fall in a viscous fluid5) Mass-spring system6) Mechanical waves7) Heat, temperature, and the body's capacity to store energy8) Linear thermal expansion9) The fundamental modes of heat transfer10) Specific and latent heat11) The Boyle-Mariotte’s Law"

# Row 18
$ws.Range("A18").Value = "Avaliação:"

# Row 19
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# Row 22
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 2, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 2, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física II, Vol. 2,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008)"
$ws.Range("C22").Value = "1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 2, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 2, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física II, Vol. 2,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008)"

# Row 23
$ws.Range("A23").Value = "Requisitos:"

# Row 24
$ws.Range("B24").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1018 -  Física I  (Requisito fraco)`n"

# Row 25
$ws.Range("B25").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"

# Fix column definitions: column A should have its own col entry (was merged with column B)
$ws.Columns.Item(1).ColumnWidth = 30.7109375